# V 0.40-B36 partial update - add new aircraft row "Kodiak 100 II Wheel" / "SWS Kodiak 100 II"
# to both the Tabelle1 (raw data) and Tabelle2 (computed/export) sheets.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Tabelle1")
$ws2 = $wb.Worksheets.Item("Tabelle2")

# ---------------------------------------------------------------------------
# Tabelle1: new row 39 - copy formatting from the most similar existing row
# (row 13, "Cessna 208B Grand Caravan EX") then overwrite with the new values.
# ---------------------------------------------------------------------------
$ws1.Range("A13:AJ13").Copy()
$ws1.Range("A39").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Set AJ39 (DispName) before A39 (ACFT) so the new shared strings are
# registered in the same order as the target file ("SWS Kodiak 100 II" then
# "Kodiak 100 II Wheel").
$ws1.Range("AJ39").Value = "SWS Kodiak 100 II"
$ws1.Range("A39").Value = "Kodiak 100 II Wheel"
$ws1.Range("B39").Value = "TB"
$ws1.Range("C39").Value = 1
$ws1.Range("D39").Value = 750
$ws1.Range("E39").Formula = "=D39/C39"
$ws1.Range("F39").Value = 0
$ws1.Range("G39").Value = 1
$ws1.Range("H39").Value = 1
$ws1.Range("I39").Value = 0
$ws1.Range("J39").Value = 0
$ws1.Range("K39").Value = 0
$ws1.Range("L39").Value = 0
$ws1.Range("M39").Value = 0
$ws1.Range("N39").Value = 1
$ws1.Range("O39").Value = 0
$ws1.Range("P39").Value = 1
$ws1.Range("Q39").Value = 1
$ws1.Range("R39").Value = 1
$ws1.Range("S39").Value = 0
$ws1.Range("T39").Value = 0
$ws1.Range("U39").Value = 0
$ws1.Range("V39").Value = 0
$ws1.Range("W39").Value = 0
$ws1.Range("X39").Value = 0
$ws1.Range("Y39").Value = 0
$ws1.Range("Z39").Value = 0
$ws1.Range("AA39").Value = 1
$ws1.Range("AB39").Value = 0
$ws1.Range("AC39").Value = 1
$ws1.Range("AD39").Value = 0
$ws1.Range("AE39").Value = 1
$ws1.Range("AF39").Value = 0
$ws1.Range("AG39").Value = 0
$ws1.Range("AH39").Value = 0
$ws1.Range("AI39").Value = 0

# ---------------------------------------------------------------------------
# Tabelle2: new row 39 - copy formatting from row 38 (the previous last row),
# then overwrite with formulas/values that mirror row 38's pattern but point
# at Tabelle1 row 39.
# ---------------------------------------------------------------------------
$ws2.Range("A38:CV38").Copy()
$ws2.Range("A39").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws2.Range("A39").Value = "|"
$ws2.Range("B39").Value = "|"
$ws2.Range("C39").Value = "|"
$ws2.Range("D39").Value = "|"
$ws2.Range("E39").Value = "|"
$ws2.Range("F39").Value = "|"
$ws2.Range("G39").Value = "|"
$ws2.Range("H39").Value = "|"
$ws2.Range("I39").Formula = "=Tabelle1!F39"
$ws2.Range("J39").Value = "|"
$ws2.Range("K39").Value = "|"
$ws2.Range("L39").Formula = "=Tabelle1!AA39"
$ws2.Range("M39").Formula = "=Tabelle1!Z39"
$ws2.Range("N39").Formula = "=Tabelle1!R39"
$ws2.Range("O39").Formula = "=Tabelle1!S39"
$ws2.Range("P39").Formula = "=Tabelle1!AC39"
$ws2.Range("Q39").Formula = "=Tabelle1!AE39"
$ws2.Range("R39").Formula = "=Tabelle1!T39"
$ws2.Range("S39").Formula = "=Tabelle1!P39"
$ws2.Range("T39").Formula = "=Tabelle1!L39"
$ws2.Range("U39").Value = "|"
$ws2.Range("V39").Value = "|"
$ws2.Range("W39").Value = "|"
$ws2.Range("X39").Value = "|"
$ws2.Range("Y39").Value = "|"
$ws2.Range("Z39").Value = "|"
$ws2.Range("AA39").Value = "|"
$ws2.Range("AB39").Value = "|"
$ws2.Range("AC39").Value = "|"
$ws2.Range("AD39").Value = "|"
$ws2.Range("AE39").Value = "|"
$ws2.Range("AF39").Value = "|"
$ws2.Range("AG39").Value = "|"
$ws2.Range("AH39").Value = "|"
$ws2.Range("AI39").Value = "|"
$ws2.Range("AJ39").Value = "|"
$ws2.Range("AK39").Value = "|"
$ws2.Range("AL39").Value = "|"
$ws2.Range("AM39").Value = "|"
$ws2.Range("AN39").Value = "|"
$ws2.Range("AO39").Value = "|"
$ws2.Range("AP39").Value = "|"
$ws2.Range("AQ39").Value = "|"
$ws2.Range("AR39").Value = "|"
$ws2.Range("AS39").Formula = "=Tabelle1!X39"
$ws2.Range("AT39").Value = "|"
$ws2.Range("AU39").Value = "|"
$ws2.Range("AV39").Value = "|"
$ws2.Range("AW39").Value = "|"
$ws2.Range("AX39").Value = "|"
$ws2.Range("AY39").Value = "|"
$ws2.Range("AZ39").Value = "|"
$ws2.Range("BA39").Value = "|"
$ws2.Range("BB39").Value = "|"
$ws2.Range("BC39").Formula = "=Tabelle1!J39"
$ws2.Range("BD39").Formula = "=Tabelle1!M39"
$ws2.Range("BE39").Formula = "=Tabelle1!AF39"
$ws2.Range("BF39").Value = "|"
$ws2.Range("BG39").Value = "|"
$ws2.Range("BH39").Value = "|"
$ws2.Range("BI39").Value = "|"
$ws2.Range("BJ39").Value = "|"
$ws2.Range("BK39").Value = "|"
$ws2.Range("BL39").Value = "|"
$ws2.Range("BM39").Value = "|"
$ws2.Range("BN39").Value = "|"
$ws2.Range("BO39").Value = "|"
$ws2.Range("BP39").Value = "|"
$ws2.Range("BQ39").Value = "|"
$ws2.Range("BR39").Value = "|"
$ws2.Range("BS39").Value = "|"
$ws2.Range("BT39").Value = "|"
$ws2.Range("BU39").Value = "|"
$ws2.Range("BV39").Value = "|"
$ws2.Range("BW39").Value = "|"
$ws2.Range("BX39").Value = "|"
$ws2.Range("BY39").Value = "|"
$ws2.Range("BZ39").Value = "|"
$ws2.Range("CA39").Value = "|"
$ws2.Range("CB39").Value = "|"
$ws2.Range("CC39").Value = "|"
$ws2.Range("CD39").Formula = "=Tabelle1!U39"
$ws2.Range("CE39").Formula = "=Tabelle1!V39"
$ws2.Range("CF39").Formula = "=Tabelle1!W39"
$ws2.Range("CG39").Formula = "=Tabelle1!Y39"
$ws2.Range("CH39").Formula = "=Tabelle1!AD39"
$ws2.Range("CI39").Formula = "=Tabelle1!N39"
$ws2.Range("CJ39").Formula = "=Tabelle1!Q39"
$ws2.Range("CK39").Formula = "=Tabelle1!K39"
$ws2.Range("CL39").Formula = "=Tabelle1!O39"
$ws2.Range("CM39").Formula = "=Tabelle1!AG39"
$ws2.Range("CN39").Value = "|"
$ws2.Range("CO39").Formula = "=Tabelle1!AH39"
$ws2.Range("CP39").Formula = "=Tabelle1!AI39"
$ws2.Range("CQ39").Value = "|"
$ws2.Range("CR39").Value = "|"
$ws2.Range("CS39").Value = "|"
$ws2.Range("CT39").Value = "|"
$ws2.Range("CU39").Value = "X"
$ws2.Range("CV39").Formula = "=Tabelle1!AJ39"

# ---------------------------------------------------------------------------
# Update the selection shown on each sheet (matches the stored sheetView).
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("B40").Select()

$ws2.Activate()
$ws2.Range("CV44").Select()

Write-Host "done"
